$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.63"
$ws.Range("E2").Value = "'0.37%"
$ws.Range("D3").Value = "'45.45"
$ws.Range("E3").Value = "'2.66%"
$ws.Range("D4").Value = "'5.576"
$ws.Range("E4").Value = "'1.53%"
$ws.Range("D5").Value = "'0.08360"
$ws.Range("E5").Value = "'4.00%"
$ws.Range("D6").Value = "'2.099"
$ws.Range("E6").Value = "'0.82%"
$ws.Range("D7").Value = "'0.9882"
$ws.Range("E7").Value = "'3.63%"
$ws.Range("E8").Value = "'-4.19%"
$ws.Range("D9").Value = "'0.1198"
$ws.Range("E9").Value = "'4.23%"
$ws.Range("D10").Value = "'0.1927"
$ws.Range("E10").Value = "'1.57%"
$ws.Range("E11").Value = "'1.05%"
$ws.Range("D12").Value = "'0.09906"
$ws.Range("E12").Value = "'-0.81%"
$ws.Range("D13").Value = "'0.04671"
$ws.Range("E13").Value = "'-2.77%"
$ws.Range("D14").Value = "'0.1059"
$ws.Range("E14").Value = "'-0.59%"
$ws.Range("D15").Value = "'0.001294"
$ws.Range("E15").Value = "'1.94%"
$ws.Range("D16").Value = "'0.005882"
$ws.Range("E16").Value = "'1.02%"
$ws.Range("D17").Value = "'3.396"
$ws.Range("E17").Value = "'0.85%"
$ws.Range("D18").Value = "'4.450"
$ws.Range("E18").Value = "'0.95%"
$ws.Range("D19").Value = "'0.3338"
$ws.Range("E19").Value = "'-3.44%"
$ws.Range("D20").Value = "'0.1373"
$ws.Range("E20").Value = "'-1.44%"
$ws.Range("D21").Value = "'0.2564"
$ws.Range("E21").Value = "'-0.69%"
$ws.Range("E22").Value = "'1.74%"
$ws.Range("D23").Value = "'0.001293"
$ws.Range("E23").Value = "'1.63%"
$ws.Range("D24").Value = "'0.004533"
$ws.Range("E24").Value = "'4.08%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'8.46%"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("D38").Value = "'0.02703"
$ws.Range("E38").Value = "'4.37%"
$ws.Range("D39").Value = "'0.05754"
$ws.Range("E39").Value = "'-1.30%"
$ws.Range("D40").Value = "'0.007907"
$ws.Range("E40").Value = "'4.31%"
$ws.Range("E41").Value = "'2.17%"
$ws.Range("D42").Value = "'0.007616"
$ws.Range("E42").Value = "'6.36%"
$ws.Range("D43").Value = "'0.002022"
$ws.Range("E43").Value = "'0.35%"
$ws.Range("D44").Value = "'0.008928"
$ws.Range("E44").Value = "'-1.65%"
$ws.Range("D45").Value = "'0.3408"
$ws.Range("D46").Value = "'0.00007048"
$ws.Range("E46").Value = "'0.80%"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("E48").Value = "'0.31%"
$ws.Range("D49").Value = "'0.003451"
$ws.Range("E49").Value = "'-2.44%"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("D51").Value = "'0.00002102"
$ws.Range("E51").Value = "'0.12%"
